$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.859.75"
$ws.Range("E2").Value = "  -2.73%  "

$ws.Range("D3").Value = "1.792.29"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5325"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3841"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07423"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.67%  "

$ws.Range("E11").Value = "  -2.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.187"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.449"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "

$ws.Range("D16").Value = "1.791.55"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.39%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06518"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.955"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "27.902.04"
$ws.Range("E23").Value = "  -2.57%  "

$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.090"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "

$ws.Range("E27").Value = "  -1.51%  "

$ws.Range("D28").Value = "1.997.62"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.315"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1092"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.23%  "

$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.647"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.500"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06935"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2202"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02270"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.037"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.375"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6100"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.175"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.04%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.682"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5694"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.911"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("E49").Value = "  +1.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06793"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.36%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000292"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +36.42%  "
